$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.576.77'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.687.93'
$ws.Range("E3").Value = '  +0.09%  '
$cell = $ws.Range("D4")
$cell.Value = "'1.006"
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.29%  '
$cell = $ws.Range("D5")
$cell.Value = "'313.76"
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$cell = $ws.Range("D6")
$cell.Value = "'1.005"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '
$cell = $ws.Range("D7")
$cell.Value = "'0.3894"
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -1.02%  '
$cell = $ws.Range("D8")
$cell.Value = "'0.4013"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = "'1.486"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '
$cell = $ws.Range("D10")
$cell.Value = "'1.006"
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +0.32%  '
$cell = $ws.Range("D11")
$cell.Value = "'52.80"
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.93%  '
$cell = $ws.Range("D12")
$cell.Value = "'0.08714"
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.93%  '
$cell = $ws.Range("D13")
$cell.Value = "'7.599"
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +5.09%  '
$cell = $ws.Range("D14")
$cell.Value = "'24.63"
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +5.20%  '
$cell = $ws.Range("D15")
$cell.Value = "'7.938"
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.25%  '
$cell = $ws.Range("D16")
$cell.Value = "'0.00001334"
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("D17").Value = '1.675.47'
$ws.Range("E17").Value = '  -1.00%  '
$cell = $ws.Range("D18")
$cell.Value = "'97.88"
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -1.60%  '
$cell = $ws.Range("D19")
$cell.Value = "'0.07095"
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.29%  '
$cell = $ws.Range("D20")
$cell.Value = "'19.64"
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.10%  '
$cell = $ws.Range("D21")
$cell.Value = "'7.241"
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +3.71%  '
$cell = $ws.Range("D22")
$cell.Value = "'1.005"
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$cell = $ws.Range("D23")
$cell.Value = "'14.15"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("D24").Value = '24.581.70'
$ws.Range("E24").Value = '  -0.30%  '
$cell = $ws.Range("D25")
$cell.Value = "'2.993"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -8.72%  '
$cell = $ws.Range("D26")
$cell.Value = "'2.346"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.67%  '
$cell = $ws.Range("D27")
$cell.Value = "'22.55"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '
$cell = $ws.Range("D28")
$cell.Value = "'161.06"
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.84%  '
$cell = $ws.Range("D29")
$cell.Value = "'8.539"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +12.83%  '
$cell = $ws.Range("D30")
$cell.Value = "'5.225"
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.87%  '
$cell = $ws.Range("D31")
$cell.Value = "'135.98"
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").Value = '1.860.44'
$ws.Range("E32").Value = '  -0.96%  '
$cell = $ws.Range("D33")
$cell.Value = "'0.08751"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.54%  '
$cell = $ws.Range("D34")
$cell.Value = "'7.449"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.26%  '
$cell = $ws.Range("D35")
$cell.Value = "'1.031"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -2.44%  '
$cell = $ws.Range("D36")
$cell.Value = "'1.976"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.95%  '
$cell = $ws.Range("D37")
$cell.Value = "'0.02895"
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +6.85%  '
$cell = $ws.Range("D38")
$cell.Value = "'0.2707"
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -0.59%  '
$cell = $ws.Range("D39")
$cell.Value = "'10.70"
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -3.86%  '
$cell = $ws.Range("D40")
$cell.Value = "'0.09086"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.74%  '
$cell = $ws.Range("D41")
$cell.Value = "'13.98"
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.48%  '
$cell = $ws.Range("D42")
$cell.Value = "'0.7729"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +1.89%  '
$cell = $ws.Range("D43")
$cell.Value = "'1.449"
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.78%  '
$cell = $ws.Range("D44")
$cell.Value = "'16.53"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +3.37%  '
$cell = $ws.Range("D45")
$cell.Value = "'0.7118"
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '
$cell = $ws.Range("D46")
$cell.Value = "'2.563"
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -1.14%  '
$cell = $ws.Range("D47")
$cell.Value = "'4.198"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.42%  '
$cell = $ws.Range("D48")
$cell.Value = "'1.006"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.41%  '
$cell = $ws.Range("D49")
$cell.Value = "'1.352"
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +3.32%  '
$cell = $ws.Range("D50")
$cell.Value = "'138.11"
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("E51").Value = '  +1.55%  '
